$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 30 (shifts existing rows 30-69 down to 31-70)
$ws.Rows(30).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(30, 1).Value = 5
$ws.Cells.Item(30, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(30, 3).Value = "Maule"
$ws.Cells.Item(30, 4).Value = 44494
$ws.Cells.Item(30, 5).Value = 7
$ws.Cells.Item(30, 6).Value = 100112013
$ws.Cells.Item(30, 7).Value = "Alcachofa"
$ws.Cells.Item(30, 8).Value = "Madrigal"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 200
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 10000
$ws.Cells.Item(30, 13).Value = 10000
$ws.Cells.Item(30, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(30, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(30, 16).Value = 250
$ws.Cells.Item(30, 17).Value = 40
$ws.Cells.Item(30, 18).Value = "Hortaliza"
